$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (cube metadata package change)
$ws.Name = "Informe-04-040001-TM"

# Drop the trailing, unused columns T:Y so the used range shrinks from A:Y to A:S
$ws.Range("T1:Y1").EntireColumn.Delete()

# Re-apply column widths for the new layout (values chosen so the stored
# OOXML width - which Excel derives from pixel-rounded ColumnWidth - lands
# on the widths used by the edited workbook)
$ws.Columns.Item(1).ColumnWidth = 26.83333333333333    # A -> 27.69
$ws.Columns.Item(2).ColumnWidth = 43.666666666666664   # B -> 44.5
$ws.Columns.Item(3).ColumnWidth = 17.833333333333332   # C -> 18.66
$ws.Columns.Item(4).ColumnWidth = 54.333333333333336   # D -> 55.2
$ws.Columns.Item(5).ColumnWidth = 33.833333333333336   # E -> 34.64
$ws.Columns.Item(6).ColumnWidth = 35.5                 # F -> 36.31
$ws.Columns.Item(7).ColumnWidth = 46.5                 # G -> 47.28
$ws.Range("H1:J1").EntireColumn.ColumnWidth = 26.83333333333333  # H:J -> 27.69
$ws.Columns.Item(11).ColumnWidth = 14.666666666666666  # K -> 15.46
$ws.Columns.Item(12).ColumnWidth = 45.666666666666664  # L -> 46.44
$ws.Columns.Item(13).ColumnWidth = 18.333333333333332  # M -> 19.19
$ws.Columns.Item(14).ColumnWidth = 19.166666666666668  # N -> 20.05
$ws.Columns.Item(15).ColumnWidth = 28.666666666666668  # O -> 29.5

# Row 6 becomes a plain empty row with the same height as rows 7/8
$ws.Rows.Item(6).RowHeight = 12.8

# A new trailing empty row 9 is appended with the same height
$ws.Rows.Item(9).RowHeight = 12.8

# Update the active selection
$ws.Range("B18").Select()
